$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 141, pushing the existing row 141 down to 142.
$ws.Rows.Item(141).Insert()

# New row 141: restore the values that used to be in row 140 before this edit
# (same Puerro / Vega Central Mapocho de Santiago record, weekly report for 2023-05-31).
$ws.Range("A141").Value = 9
$ws.Range("B141").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C141").Value = "Metropolitana"
$ws.Range("D141").Value = 45077
$ws.Range("D141").NumberFormat = $ws.Range("D140").NumberFormat
$ws.Range("E141").Value = 13
$ws.Range("F141").Value = 100112005
$ws.Range("G141").Value = "Puerro"
$ws.Range("H141").Value = "Sin especificar"
$ws.Range("I141").Value = "Primera"
$ws.Range("J141").Value = 70
$ws.Range("K141").Value = 8000
$ws.Range("L141").Value = 8000
$ws.Range("M141").Value = 8000
$ws.Range("N141").Value = "$/paquete 20 unidades"
$ws.Range("O141").Value = "Provincia de Chacabuco"
$ws.Range("P141").Value = 400
$ws.Range("Q141").Value = 20
$ws.Range("R141").Value = "Hortaliza"

# Row 139: new weekly record (2023-11-09), volume dropped to 160.
$ws.Range("D139").Value = 45239
$ws.Range("J139").Value = 160

# Row 140: shifted back to the 2021-05-19 record's figures.
$ws.Range("D140").Value = 44335
$ws.Range("J140").Value = 250
$ws.Range("K140").Value = 7000
$ws.Range("L140").Value = 8000
$ws.Range("M140").Value = 7500
$ws.Range("P140").Value = 375
